$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before column E (old E..J shift right to I..N),
# making room for the new MySQL connection-info columns.
$ws.Range("E1:H1").EntireColumn.Insert()

# The inserted columns pick up E1's old formatting (style s="2"); clear it so
# the new header/data cells fall back to the sheet's default (unstyled) look.
$ws.Range("E1:H2").ClearFormats()

# --- Row 1 headers (new MySQL columns) ---
$ws.Range("E1").Value = "SqlIP"
$ws.Range("F1").Value = "SqlPort"
$ws.Range("G1").Value = "SqlName"
$ws.Range("H1").Value = "SqlUser"
$ws.Range("I1").Value = "SqlPwd"

# --- Row 2 MySQL connection values ---
$ws.Range("E2").Value = "192.168.0.24"
$ws.Range("F2").Value = 3306
$ws.Range("G2").Value = "app_test"
$ws.Range("H2").Value = "root"
# I2 already holds the old E2 value (123456) after the column insert shifted it.

# --- Column widths (character units, as stored in the sheet). The host's
# column width is quantized to 1/7-character pixel steps on save, so these
# inputs are chosen to round-trip to the closest possible value to the
# target stored width (exact fractional eighths aren't reachable via
# ColumnWidth's pixel rounding).
$ws.Range("A1").EntireColumn.ColumnWidth = 11.929
$ws.Range("D1").EntireColumn.ColumnWidth = 4.643
$ws.Range("E1").EntireColumn.ColumnWidth = 13.072
$ws.Range("F1").EntireColumn.ColumnWidth = 7.643
$ws.Range("G1").EntireColumn.ColumnWidth = 7.643
$ws.Range("H1").EntireColumn.ColumnWidth = 7.643
$ws.Range("I1").EntireColumn.ColumnWidth = 6.643

# --- Selection moves to G9 ---
$ws.Range("G9").Select()

# --- Data validation: extend the "TRUE,FALSE" list rule to the new header
# cells, and keep the blank-allowed placeholder on the (now shifted) J1 cell.
$ws.Range("J2:J1048576").Validation.Delete()
$ws.Range("J2:J11").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws.Range("J13:J1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws.Range("I1").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws.Range("E1:F1").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws.Range("H1").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
